# Edit: elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 6 new rows starting at row 18 -----------------------------
# Before: rows 16-18 hold the 3 data rows (row 18 has the "closing" bottom
# border style). Inserting at row 18 pushes the old row18 (and everything
# below it, including the signature rows 23-24) down by 6 rows, so the
# "closing" styled row ends up at row 24, and the signature block ends up
# at rows 29-30, matching the target layout (B2:J30).
$ws.Rows("18:23").Insert()

# Copy the formatting (borders/fonts/number formats) of the "middle" data
# row (row 17) into the freshly inserted rows 18:23 so they look like
# normal (non-closing) table rows.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J23").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Update summary fields ---------------------------------------------
$ws.Range("E11").Value2 = 483293   # VALOR MORA (total)
$ws.Range("C13").Value2 = 6        # Cant. Trabajadores
$ws.Range("F13").Value2 = 7        # Cant. Periodos

# --- 3. Rewrite the data table (rows 16-24) --------------------------------
$data = @(
    @("CC", "73169933",   "ADALBERTO RAFAEL MARTINEZ NOVOA", "2507", 56940,  986000),
    @("CC", "73101200",   "VIRGILIO PIANETA GONZALEZ",        "2507", 120000, 3000000),
    @("CC", "1143380904", "MANUEL SALVADOR VILLA CARO",       "2408", 52000,  908526),
    @("CC", "1143380904", "MANUEL SALVADOR VILLA CARO",       "2407", 52000,  908526),
    @("CC", "1143380904", "MANUEL SALVADOR VILLA CARO",       "2406", 52000,  908526),
    @("CC", "1143380904", "MANUEL SALVADOR VILLA CARO",       "2405", 39866,  908526),
    @("CC", "1047419802", "IVAN ENRIQUE CANOLES GONZALEZ",    "2507", 56940,  1423500),
    @("CC", "1002388492", "JHON FREDIS HEIBAN ROMERO DIAZ",   "2409", 52000,  1300000),
    @("CC", "73009379",   "GREGORIO CARRILLO RODRIGUEZ",      "2311", 1547,   1160000)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value2 = $rec[0]
    $ws.Range("C$row").Value2 = $rec[1]
    $ws.Range("D$row").Value2 = $rec[2]
    $ws.Range("E$row").Value2 = $rec[3]
    $ws.Range("F$row").Value2 = $rec[4]
    $ws.Range("G$row").Value2 = $rec[5]
    $row = $row + 1
}
